# Apply the "generated excel file from csv files as inputs" update to the
# "stories" worksheet: add a missing J1 header column and refresh the
# estimated/consumed/pending effort figures (and a couple of derived
# percentages) for rows 2-15 based on the newly generated CSV data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stories")

# --- New header cell J1 (previously missing) ---------------------------
$ws.Range("J1").Value = "ESTIMATED`nEFFORT`n(IN HRS))"

# --- Row 2 ---------------------------------------------------------------
$ws.Range("J2").Value = 16
$ws.Range("K2").Value = 0
$ws.Range("M2").Value = 0

# --- Row 3 ---------------------------------------------------------------
$ws.Range("J3").Value = 16
$ws.Range("K3").Value = 0
$ws.Range("M3").Value = 0

# --- Row 4 ---------------------------------------------------------------
$ws.Range("J4").Value = 16
$ws.Range("K4").Value = 0
$ws.Range("M4").Value = 0

# --- Row 5 ---------------------------------------------------------------
$ws.Range("J5").Value = 16
$ws.Range("K5").Value = 43200
$ws.Range("M5").Value = 0
$ws.Range("P5").Value = 100

# --- Row 6 ---------------------------------------------------------------
$ws.Range("J6").Value = 16
$ws.Range("K6").Value = 0
$ws.Range("M6").Value = 0

# --- Row 7 ---------------------------------------------------------------
$ws.Range("J7").Value = 16
$ws.Range("K7").Value = 0
$ws.Range("M7").Value = 0

# --- Row 8 ---------------------------------------------------------------
$ws.Range("J8").Value = 24
$ws.Range("K8").Value = 82800
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 23
$ws.Range("P8").Value = 4

# --- Row 9 ---------------------------------------------------------------
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 16

# --- Row 10 --------------------------------------------------------------
$ws.Range("C10").Value = "PCS-2122"
$ws.Range("J10").Value = 32
$ws.Range("K10").Value = 28800
$ws.Range("L10").Value = 24
$ws.Range("M10").Value = 8
$ws.Range("P10").Value = 75

# --- Row 11 --------------------------------------------------------------
$ws.Range("J11").Value = 16
$ws.Range("K11").Value = 0
$ws.Range("M11").Value = 0

# --- Row 12 --------------------------------------------------------------
$ws.Range("J12").Value = 8
$ws.Range("K12").Value = 0
$ws.Range("M12").Value = 0

# --- Row 13 --------------------------------------------------------------
$ws.Range("J13").Value = 16
$ws.Range("K13").Value = 0
$ws.Range("M13").Value = 0

# --- Row 14 --------------------------------------------------------------
$ws.Range("J14").Value = 16
$ws.Range("K14").Value = 0
$ws.Range("M14").Value = 0

# --- Row 15 --------------------------------------------------------------
$ws.Range("J15").Value = 8
$ws.Range("K15").Value = 14400
$ws.Range("L15").Value = 4
$ws.Range("P15").Value = 50
